# Append the latest Adafruit IO reading as a new row (row 8) to the sheet.
# Mirrors the existing rows: every column is stored as plain text, including
# numeric-looking values like the "Value" column, so a leading apostrophe is
# used to force text entry for that cell (then the quote-prefix style that
# Excel applies is cleared again so no stray cell style is introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
